$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.385.35'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '1.875.75'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.51'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  -2.86%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06511'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').Value = '1.871.82'
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07486'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.55'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '88.35'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6603'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Value = '30.345.24'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007581'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('D20').Value = '2.115.27'
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.297'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.39%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '219.85'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +14.17%  '
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.329'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '167.73'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.41'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.959'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09375'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.302'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.017'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05019'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.204'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.54%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7417'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.12%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.709'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01821'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.615'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.060'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9034'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '106.40'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.854'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.4264'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.73%  '
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.405'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.30%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '64.34'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('E47').Value = '  -7.40%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.472'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -6.51%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.890'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '33.66'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.85%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05630'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.59%  '
